$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 405.2381
$ws.Range("I19").Value = 345.85715
$ws.Range("J19").Value = 434.92856
$ws.Range("K19").Value = 345.85715
$ws.Range("L19").Value = 434.92856
$ws.Range("M19").Value = -170.85715
$ws.Range("N19").Value = -784.9285600000001

$ws.Range("H74").Value = 14998
$ws.Range("I74").Value = 70003
$ws.Range("K74").Value = 70003
$ws.Range("M74").Value = -69067

$ws.Range("H77").Value = 14998
$ws.Range("I77").Value = 70003
$ws.Range("K77").Value = 350015
$ws.Range("M77").Value = -345335

$ws.Range("H113").Value = 26199.4
$ws.Range("I113").Value = 23249.25
$ws.Range("K113").Value = 23249.25
$ws.Range("M113").Value = -19995.25

$ws.Range("H129").Value = 2356.4285

$ws.Range("H132").Value = 8712.25
$ws.Range("I132").Value = 9458
$ws.Range("K132").Value = 28374
$ws.Range("M132").Value = -25844

$ws.Range("H137").Value = 9480.444
$ws.Range("I137").Value = 12104.053
$ws.Range("K137").Value = 36312.159
$ws.Range("M137").Value = -33762.159

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 6637.5
$ws.Range("I8").Value = 327
$ws.Range("K8").Value = 327
$ws.Range("M8").Value = -183

$ws.Range("H12").Value = 199.5
$ws.Range("J12").Value = 199
$ws.Range("L12").Value = 199
$ws.Range("N12").Value = -545

$ws.Range("H32").Value = 5084.1465
$ws.Range("I32").Value = 4979.6025
$ws.Range("J32").Value = 8900
$ws.Range("K32").Value = 4979.6025
$ws.Range("L32").Value = 8900
$ws.Range("M32").Value = -4692.6025
$ws.Range("N32").Value = -9474

$ws.Range("H44").Value = 49999.75
$ws.Range("J44").Value = 49999.75
$ws.Range("L44").Value = 49999.75
$ws.Range("N44").Value = -50975.75

$ws.Range("J63").Value = 2000
$ws.Range("L63").Value = 2000
$ws.Range("N63").Value = -3372

$ws.Range("J66").Value = 2000
$ws.Range("L66").Value = 10000
$ws.Range("N66").Value = -16864

$ws.Range("H88").Value = 71429790
$ws.Range("I88").Value = 928.5
$ws.Range("K88").Value = 928.5
$ws.Range("M88").Value = -522.5

$ws.Range("H91").Value = 71429790
$ws.Range("I91").Value = 928.5
$ws.Range("K91").Value = 928.5
$ws.Range("M91").Value = 475.5

$ws.Range("H107").Value = 10000
$ws.Range("I107").Value = 10000
$ws.Range("K107").Value = 10000
$ws.Range("M107").Value = -6160

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8663
$ws.Range("I20").Value = 5808
$ws.Range("J20").Value = 10566.333
$ws.Range("K20").Value = 5808
$ws.Range("L20").Value = 10566.333
$ws.Range("M20").Value = -5561
$ws.Range("N20").Value = -11060.333

$ws.Range("H105").Value = 207600
$ws.Range("J105").Value = 9500
$ws.Range("L105").Value = 9500
$ws.Range("N105").Value = -12994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 508
$ws.Range("I19").Value = 508
$ws.Range("K19").Value = 508
$ws.Range("M19").Value = -338

$ws.Range("H24").Value = 508
$ws.Range("I24").Value = 508
$ws.Range("K24").Value = 508
$ws.Range("M24").Value = -338

$ws.Range("H99").Value = 10420822
$ws.Range("I99").Value = 13891762
$ws.Range("J99").Value = 8000
$ws.Range("K99").Value = 13891762
$ws.Range("L99").Value = 8000
$ws.Range("M99").Value = -13890264
$ws.Range("N99").Value = -10996

$ws.Range("H105").Value = 177877.5
$ws.Range("I105").Value = 211753
$ws.Range("K105").Value = 211753
$ws.Range("M105").Value = -210006

$ws.Range("H126").Value = 10420822
$ws.Range("I126").Value = 13891762
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 41675286
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -41672816
$ws.Range("N126").Value = -28940

$ws.Range("H132").Value = 2665.6667
$ws.Range("I132").Value = 2665.6667
$ws.Range("K132").Value = 7997.000100000001
$ws.Range("M132").Value = -5467.000100000001

$ws.Range("H134").Value = 3484.5293
$ws.Range("I134").Value = 3853.3333
$ws.Range("J134").Value = 2599.4
$ws.Range("K134").Value = 11559.9999
$ws.Range("L134").Value = 7798.200000000001
$ws.Range("M134").Value = -9024.999899999999
$ws.Range("N134").Value = -12868.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 345.5
$ws.Range("I10").Value = 148.33333
$ws.Range("J10").Value = 463.8
$ws.Range("K10").Value = 444.99999
$ws.Range("L10").Value = 1391.4
$ws.Range("M10").Value = -305.99999
$ws.Range("N10").Value = -1669.4

$ws.Range("H25").Value = 2240.8948
$ws.Range("I25").Value = 1750
$ws.Range("J25").Value = 2298.647
$ws.Range("K25").Value = 5250
$ws.Range("L25").Value = 6895.941
$ws.Range("M25").Value = -5081
$ws.Range("N25").Value = -7233.941

$ws.Range("H30").Value = 2240.8948
$ws.Range("I30").Value = 1750
$ws.Range("J30").Value = 2298.647
$ws.Range("K30").Value = 5250
$ws.Range("L30").Value = 6895.941
$ws.Range("M30").Value = -5148
$ws.Range("N30").Value = -7099.941

$ws.Range("H105").Value = 9999.924999999999
$ws.Range("J105").Value = 9999.924999999999
$ws.Range("L105").Value = 29999.775
$ws.Range("N105").Value = -35241.77499999999

$ws.Range("H122").Value = 5696.483
$ws.Range("I122").Value = 190
$ws.Range("J122").Value = 6104.3706
$ws.Range("K122").Value = 1710
$ws.Range("L122").Value = 54939.3354
$ws.Range("M122").Value = 740
$ws.Range("N122").Value = -59839.3354

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 9619.559999999999
$ws.Range("I122").Value = 6026.0527
$ws.Range("J122").Value = 20999
$ws.Range("K122").Value = 18078.1581
$ws.Range("L122").Value = 62997
$ws.Range("M122").Value = -15628.1581
$ws.Range("N122").Value = -67897

$ws.Range("H132").Value = 3258.652
$ws.Range("I132").Value = 3258.652
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9775.956
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7245.956
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 37799.6
$ws.Range("J136").Value = 37799.6
$ws.Range("L136").Value = 113398.8
$ws.Range("N136").Value = -118498.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4818.875
$ws.Range("I82").Value = 4925.1665
$ws.Range("J82").Value = 4500
$ws.Range("K82").Value = 4925.1665
$ws.Range("L82").Value = 4500
$ws.Range("M82").Value = -4564.1665
$ws.Range("N82").Value = -5222

$ws.Range("H85").Value = 4818.875
$ws.Range("I85").Value = 4925.1665
$ws.Range("J85").Value = 4500
$ws.Range("K85").Value = 4925.1665
$ws.Range("L85").Value = 4500
$ws.Range("M85").Value = -3677.1665
$ws.Range("N85").Value = -6996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 27797.562
$ws.Range("I126").Value = 30340.428
$ws.Range("K126").Value = 91021.284
$ws.Range("M126").Value = -88551.284

$ws.Range("H137").Value = 59497.75
$ws.Range("J137").Value = 59497.75
$ws.Range("L137").Value = 59497.75
$ws.Range("N137").Value = -69697.75
